$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting the existing row 57 (and all
# rows below it, down through the former row 97) down by one row. This
# mirrors Excel's "Insert Sheet Rows" behaviour (equivalent to
# xlShiftDown), which is exactly what the target diff shows: the old
# rows 57..97 reappear unchanged as rows 58..98, and a brand-new row of
# data is written into the now-empty row 57.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly price entry.
$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).Value = 45062
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = 100112043
$ws.Cells.Item(57, 7).Value = "Pepino dulce"
$ws.Cells.Item(57, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 80
$ws.Cells.Item(57, 11).Value = 20000
$ws.Cells.Item(57, 12).Value = 21000
$ws.Cells.Item(57, 13).Value = 20500
$ws.Cells.Item(57, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(57, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(57, 16).Value = 1139
$ws.Cells.Item(57, 17).Value = 18
$ws.Cells.Item(57, 18).Value = "Hortaliza"
